# Applies the "Project Deliverable" content edits described by the diff:
#   1. Rewrite the MarketPl.ai intro paragraph (3 runs -> 1 run, new copy).
#   2. Rewrite the "You will create..." paragraph with the new feature summary.
#   3. Merge "Key prompts or interactions" + " (2-3 examples are enough)" runs.
#   4. Merge "Include at least " / "two screenshots..." / "." runs.
#   5. Merge "Total size: " / "typically" / " 3-6 pages" runs (drops proofErr tags).

$d = $word.ActiveDocument

# --- 1. MarketPl.ai intro paragraph -----------------------------------
# Scope the Find to the paragraph's own Range so the identical "MarketPl.ai"
# text in the title (paragraph 1) is left untouched.
$introPara = $d.Paragraphs.Item(3)
$introRange = $introPara.Range
$introRange.Find.Execute(
    "MarketPl.ai helps users evaluate how different investment strategies would have performed using historical market data. Unlike live trading platforms or brokerage tools, this application focuses on offline analysis, clarity, and educational insight - allowing users to explore how strategies such as Buy & Hold, Moving Average Crossover, or Periodic Investing behave over time under real historical conditions.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "MarketPl.ai is a web application for exploring how investment strategies would have performed using historical stock and ETF price data. It is designed for analysis and learning: users select a symbol and time period, run simulations under consistent assumptions, and review outcomes through clear charts, metrics, and trade history.",
    2
)

# --- 2. "You will create..." -> new feature summary --------------------
$d.Content.Find.Execute(
    "You will create a modern desktop or web application that enables users to import historical price data from CSV files, configure investment strategies, simulate portfolio evolution, and analyze return and risk metrics through clear visualizations and comparisons.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Users start by importing historical price data from CSV files and monitoring import progress and status. Once data is available, they can browse symbols, filter by date ranges, and inspect price series in a data explorer. Backtesting features allow users to choose a strategy (e.g., Buy and Hold, Dollar-Cost Averaging, Moving Average Crossover, RSI, Bollinger Bands), configure strategy parameters, and run simulations for a single symbol or a small weighted portfolio. Results include an equity curve, key performance and risk metrics, and a trade log, with the option to compare multiple strategies side-by-side on the same dataset.",
    2
)

# --- 3. "Key prompts or interactions" + " (2-3 examples are enough)" ---
$d.Content.Find.Execute(
    "Key prompts or interactions (2–3 examples are enough)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Key prompts or interactions (2–3 examples are enough)",
    2
)

# --- 4. "Include at least " / "two screenshots..." / "." ---------------
$d.Content.Find.Execute(
    "Include at least two screenshots of the functioning system.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Include at least two screenshots of the functioning system.",
    2
)

# --- 5. "Total size: " / "typically" / " 3-6 pages" ---------------------
$d.Content.Find.Execute(
    "Total size: typically 3-6 pages",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Total size: typically 3-6 pages",
    2
)
